$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the user data so B4 and B5 both hold the correct "secret_sauce" value
$ws.Range("B4").Value = "secret_sauce"
$ws.Range("B5").Value = "secret_sauce"

# Move/select B5 as the active cell, matching the updated selection in the sheet view
$ws.Range("B5").Select()
